$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 241, shifting existing rows 241:343 down to 242:344
$ws.Rows.Item(241).Insert()

# Populate the newly inserted row 241 with the new record's data
$ws.Cells.Item(241, 1).Value = 5
$ws.Cells.Item(241, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(241, 3).Value = "Maule"
$ws.Cells.Item(241, 4).Value = 44726
$ws.Cells.Item(241, 5).Value = 7
$ws.Cells.Item(241, 6).Value = 100112023
$ws.Cells.Item(241, 7).Value = "Brócoli"
$ws.Cells.Item(241, 8).Value = "Sin especificar"
$ws.Cells.Item(241, 9).Value = "Primera"
$ws.Cells.Item(241, 10).Value = 5000
$ws.Cells.Item(241, 11).Value = 800
$ws.Cells.Item(241, 12).Value = 800
$ws.Cells.Item(241, 13).Value = 800
$ws.Cells.Item(241, 14).Value = "$/unidad"
$ws.Cells.Item(241, 15).Value = "Región del Maule"
$ws.Cells.Item(241, 16).Value = 800
$ws.Cells.Item(241, 17).Value = 1
$ws.Cells.Item(241, 18).Value = "Hortaliza"
